$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "16-03-25"
$ws.Range("B3").Value = "Viju"
$ws.Range("D3").Value = "16-03-25"
$ws.Range("A4").Value = "www"
$ws.Range("C4").Value = "nds.nitin@gmail.com "
$ws.Range("D4").Value = "16-03-25"
